# pre-processing notebook for ammonia data
# Add a "lhv" sheet (lower heating values) alongside the existing "hhv"
# sheet, built by duplicating "hhv" and swapping in the LHV figures.

$wb = $excel.ActiveWorkbook
$hhv = $wb.Worksheets.Item("hhv")

# Capture the selection on "hhv" while it is still the active sheet.
$hhv.Range("A3:A16").Select() | Out-Null

# Duplicate "hhv" to the end of the workbook; the copy becomes "lhv".
$hhv.Copy($null, $hhv)
$lhv = $wb.Worksheets.Item($wb.Worksheets.Count)
$lhv.Name = "lhv"

# Drop the duplicated "Natural Gas" row and the trailing "Average" row
# (highest row index first so the remaining row numbers don't shift).
$lhv.Rows.Item(15).Delete() | Out-Null
$lhv.Rows.Item(9).Delete() | Out-Null

# Swap the HHV figures for the corresponding LHV figures.
$lhv.Range("B3").Value = 47.141
$lhv.Range("B4").Value = 42.791
$lhv.Range("B5").Value = 43.448
$lhv.Range("B6").Value = 39.466
$lhv.Range("B7").Value = 42.612
$lhv.Range("B8").Value = 22.732
$lhv.Range("B9").Value = 26.122
$lhv.Range("B10").Value = 19.551
$lhv.Range("B11").Value = 37.528
$lhv.Range("B12").Value = 26.952
$lhv.Range("B13").Value = 16.37
$lhv.Range("B14").Value = 46.898

# Leave the new sheet active with its own selection, matching the
# author's last on-screen state.
$lhv.Range("D10").Select() | Out-Null
